$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "want to go" counts (column F)
$ws1 = $wb.Worksheets("展览")
$ws1.Range("F2").Value = 3377
$ws1.Range("F5").Value = 6958
$ws1.Range("F6").Value = 2383
$ws1.Range("F8").Value = 103
$ws1.Range("F12").Value = 32
$ws1.Range("F13").Value = 172
$ws1.Range("F14").Value = 565

# Sheet "演出" (Performance) - update "want to go" counts (column F)
$ws2 = $wb.Worksheets("演出")
$ws2.Range("F2").Value = 22

# Sheet "全部类型" (All types) - update "want to go" counts (column F)
$ws4 = $wb.Worksheets("全部类型")
$ws4.Range("F2").Value = 3377
$ws4.Range("F3").Value = 22
$ws4.Range("F6").Value = 6958
$ws4.Range("F7").Value = 2383
$ws4.Range("F9").Value = 103
$ws4.Range("F13").Value = 32
$ws4.Range("F14").Value = 172
$ws4.Range("F15").Value = 565
